# Add a "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: same style as the other header cells (bold/centered/bordered) -
# copy the format from G1 (the "sum" header) so we reuse the existing style index
# instead of creating a brand-new one.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Save values per row (1 if the pitcher's outing counted as a save, else 0)
$saveValues = @{
    2 = 0; 3 = 0; 4 = 0; 5 = 0; 6 = 0; 7 = 0; 8 = 0; 9 = 0; 10 = 1;
    11 = 0; 12 = 0; 13 = 0; 14 = 0; 15 = 0; 16 = 1; 17 = 0; 18 = 0;
    19 = 0; 20 = 0; 21 = 0; 22 = 0; 23 = 0; 24 = 1; 25 = 1; 26 = 0;
    27 = 1; 28 = 1; 29 = 0; 30 = 0; 31 = 0; 32 = 0; 33 = 0; 34 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
